$wb = $excel.ActiveWorkbook

# Rename the existing "Taxa" sheet to "GBIFTaxa"
$gbifTaxa = $wb.Worksheets.Item("Taxa")
$gbifTaxa.Name = "GBIFTaxa"

# Update GBIFTaxa sheet's current selection (F17 -> E29)
[void]$gbifTaxa.Range("E29").Select()

# Insert a new "NCBITaxa" sheet immediately before the "Data" sheet
$dataSheet = $wb.Worksheets.Item("Data")
$ncbiTaxa = $wb.Worksheets.Add($dataSheet)
$ncbiTaxa.Name = "NCBITaxa"

# Populate the header row for the new NCBITaxa sheet
$ncbiTaxa.Range("A1").Value = "Name"
$ncbiTaxa.Range("B1").Value = "NCBI ID"
$ncbiTaxa.Range("C1").Value = "Add taxonomic ranks here"
$ncbiTaxa.Range("D1").Value = "Comments"

# Match column C width to the authored template (~22.66 chars, best-fit)
$ncbiTaxa.Columns.Item(3).ColumnWidth = 21.8307291666667

# Set the active cell/selection on the new sheet and make it the active tab
[void]$ncbiTaxa.Range("F7").Select()
[void]$ncbiTaxa.Activate()
